# Auto-generated: apply scheduled-runner market-data refresh to Sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 16
$ws.Range("I8").Value = 17.666666
$ws.Range("J8").Value = 15
$ws.Range("K8").Value = 52.999998
$ws.Range("L8").Value = 45
$ws.Range("M8").Value = 86.00000199999999
$ws.Range("N8").Value = -323

$ws.Range("H28").Value = 1409
$ws.Range("I28").Value = 701.6667
$ws.Range("K28").Value = 701.6667
$ws.Range("M28").Value = -216.6667

$ws.Range("H41").Value = 244.125
$ws.Range("I41").Value = 244.125
$ws.Range("K41").Value = 244.125
$ws.Range("M41").Value = 195.875

$ws.Range("H53").Value = 298.9
$ws.Range("I53").Value = 384.2857
$ws.Range("J53").Value = 99.666664
$ws.Range("K53").Value = 384.2857
$ws.Range("L53").Value = 99.666664
$ws.Range("M53").Value = 252.7143
$ws.Range("N53").Value = -1373.666664

$ws.Range("H74").Value = 3691.5
$ws.Range("I74").Value = 3691.5
$ws.Range("K74").Value = 3691.5
$ws.Range("M74").Value = -2755.5

$ws.Range("H77").Value = 3691.5
$ws.Range("I77").Value = 3691.5
$ws.Range("K77").Value = 18457.5
$ws.Range("M77").Value = -13777.5

$ws.Range("H100").Value = 4374.5
$ws.Range("I100").Value = 4374.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4374.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3833.5
$ws.Range("N100").ClearContents()

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 2033.8334
$ws.Range("I132").Value = 2061.3914
$ws.Range("K132").Value = 6184.174199999999
$ws.Range("M132").Value = -3654.174199999999

$ws.Range("H138").Value = 3641
$ws.Range("J138").Value = 3881.8333
$ws.Range("L138").Value = 11645.4999
$ws.Range("N138").Value = -21925.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1023
$ws.Range("J2").Value = 1056.6
$ws.Range("L2").Value = 1056.6
$ws.Range("N2").Value = -1282.6

$ws.Range("H61").Value = 4836
$ws.Range("I61").Value = 4836
$ws.Range("K61").Value = 4836
$ws.Range("M61").Value = -4624

$ws.Range("H74").Value = 3187.25
$ws.Range("I74").Value = 3252.4546
$ws.Range("K74").Value = 3252.4546
$ws.Range("M74").Value = -2378.4546

$ws.Range("H77").Value = 3187.25
$ws.Range("I77").Value = 3252.4546
$ws.Range("K77").Value = 16262.273
$ws.Range("M77").Value = -11894.273

$ws.Range("H116").Value = 1023
$ws.Range("J116").Value = 1056.6
$ws.Range("L116").Value = 1056.6
$ws.Range("N116").Value = -5644.6

$ws.Range("H132").Value = 1370.0741
$ws.Range("I132").Value = 1249.6818
$ws.Range("K132").Value = 3749.0454
$ws.Range("M132").Value = -1219.0454

$ws.Range("H136").Value = 4836
$ws.Range("I136").Value = 4836
$ws.Range("K136").Value = 14508
$ws.Range("M136").Value = -11958

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1023
$ws.Range("J3").Value = 1056.6
$ws.Range("L3").Value = 1056.6
$ws.Range("N3").Value = -1284.6

$ws.Range("H20").Value = 4664.6665
$ws.Range("I20").Value = 3698.8
$ws.Range("K20").Value = 3698.8
$ws.Range("M20").Value = -3451.8

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

$ws.Range("H134").Value = 2864.3914
$ws.Range("I134").Value = 3117.4211
$ws.Range("J134").Value = 1662.5
$ws.Range("K134").Value = 9352.263300000001
$ws.Range("L134").Value = 4987.5
$ws.Range("M134").Value = -6817.263300000001
$ws.Range("N134").Value = -10057.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 396.25
$ws.Range("I22").Value = 392.5
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 392.5
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -42.5
$ws.Range("N22").Value = -1100

$ws.Range("H105").Value = 1308.5
$ws.Range("I105").Value = 1281.1428
$ws.Range("K105").Value = 1281.1428
$ws.Range("M105").Value = 465.8571999999999

$ws.Range("H107").Value = 947.7273
$ws.Range("I107").Value = 740.75
$ws.Range("K107").Value = 740.75
$ws.Range("M107").Value = 1179.25

$ws.Range("H134").Value = 3000
$ws.Range("J134").Value = 3000
$ws.Range("L134").Value = 9000
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1166.6666
$ws.Range("J63").Value = 1500
$ws.Range("L63").Value = 4500
$ws.Range("N63").Value = -5998

$ws.Range("H66").Value = 1166.6666
$ws.Range("J66").Value = 1500
$ws.Range("L66").Value = 13500
$ws.Range("N66").Value = -20988

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3593.75
$ws.Range("I113").Value = 3958.3333
$ws.Range("K113").Value = 3958.3333
$ws.Range("M113").Value = -1788.3333

$ws.Range("H132").Value = 2076
$ws.Range("I132").Value = 2076
$ws.Range("K132").Value = 6228
$ws.Range("M132").Value = -3698

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3547.7827
$ws.Range("I46").Value = 2888.889
$ws.Range("J46").Value = 3971.3572
$ws.Range("K46").Value = 2888.889
$ws.Range("L46").Value = 3971.3572
$ws.Range("M46").Value = -2700.889
$ws.Range("N46").Value = -4347.3572

$ws.Range("H132").Value = 5319.8184
$ws.Range("I132").Value = 3900.3
$ws.Range("J132").Value = 6502.75
$ws.Range("K132").Value = 11700.9
$ws.Range("L132").Value = 19508.25
$ws.Range("M132").Value = -9170.900000000001
$ws.Range("N132").Value = -24568.25

$ws.Range("H136").Value = 3018.889
$ws.Range("I136").Value = 2896.875
$ws.Range("J136").Value = 3995
$ws.Range("K136").Value = 8690.625
$ws.Range("L136").Value = 11985
$ws.Range("M136").Value = -6140.625
$ws.Range("N136").Value = -17085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 7999.4
$ws.Range("I2").Value = 5000
$ws.Range("K2").Value = 5000
$ws.Range("M2").Value = -4888

$ws.Range("H22").Value = 4000
$ws.Range("J22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("N22").Value = -4586

$ws.Range("H100").Value = 2000
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H132").Value = 1317.3846
$ws.Range("I132").Value = 959.1667
$ws.Range("J132").Value = 1624.4286
$ws.Range("K132").Value = 2877.5001
$ws.Range("L132").Value = 4873.2858
$ws.Range("M132").Value = -347.5001000000002
$ws.Range("N132").Value = -9933.2858

$ws.Range("H136").Value = 6157.8945
$ws.Range("I136").Value = 6301.2144
$ws.Range("J136").Value = 5756.6
$ws.Range("K136").Value = 18903.6432
$ws.Range("L136").Value = 17269.8
$ws.Range("M136").Value = -16353.6432
$ws.Range("N136").Value = -22369.8
